# Added section 2 and 3
#
# The document ends with a trailing list-item paragraph (style
# "ListParagraph", same numbering as the items above it) that has no
# text in it yet. Give that paragraph the text "Added section 2", then
# add a brand-new paragraph right after it - same list formatting -
# with the text "Added section 3".

$d = $word.ActiveDocument

# The last paragraph in the document is the empty bullet waiting for
# content.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "Added section 2"

# Re-fetch the (now re-seated) last paragraph and append a sibling
# paragraph after it; InsertParagraphAfter carries over the paragraph
# formatting (style/numbering) of the paragraph it splits off from.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Added section 3"
